$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 15 and 16 swap positions: WrappedEther and Chainlink exchange ranking
$ws.Range("B15").Value = "Chainlink"
$ws.Range("C15").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D15").Value = "'7.232"
$ws.Range("E15").Value = "  -1.94%  "

$ws.Range("B16").Value = "WrappedEther"
$ws.Range("C16").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D16").Value = "1.787.69"
$ws.Range("E16").Value = "  -1.09%  "

# Price and Volume(1h) updates for remaining rows
$ws.Range("D2").Value = "28.109.66"
$ws.Range("E2").Value = "  -0.81%  "
$ws.Range("D3").Value = "1.794.59"
$ws.Range("E3").Value = "  +0.14%  "
$ws.Range("D4").Value = "'1.000"
$ws.Range("E4").Value = "  -0.52%  "
$ws.Range("D5").Value = "'317.28"
$ws.Range("E5").Value = "  +0.86%  "
$ws.Range("D6").Value = "'1.000"
$ws.Range("E6").Value = "  -0.24%  "
$ws.Range("D7").Value = "'0.5418"
$ws.Range("E7").Value = "  +0.56%  "
$ws.Range("D8").Value = "'0.3774"
$ws.Range("E8").Value = "  -1.32%  "
$ws.Range("D9").Value = "'0.07446"
$ws.Range("E9").Value = "  -1.45%  "
$ws.Range("D10").Value = "'41.68"
$ws.Range("E10").Value = "  -1.85%  "
$ws.Range("D11").Value = "'1.095"
$ws.Range("E11").Value = "  -2.13%  "
$ws.Range("D12").Value = "'1.000"
$ws.Range("E12").Value = "  -0.60%  "
$ws.Range("D13").Value = "'20.55"
$ws.Range("E13").Value = "  -2.24%  "
$ws.Range("D14").Value = "'6.116"
$ws.Range("E14").Value = "  -0.53%  "
$ws.Range("D17").Value = "'89.17"
$ws.Range("E17").Value = "  -2.60%  "
$ws.Range("D18").Value = "'0.00001058"
$ws.Range("E18").Value = "  -0.97%  "
$ws.Range("D19").Value = "'0.06482"
$ws.Range("E19").Value = "  +0.74%  "
$ws.Range("D20").Value = "'1.0000"
$ws.Range("E20").Value = "  -0.19%  "
$ws.Range("D21").Value = "'17.30"
$ws.Range("E21").Value = "  +0.23%  "
$ws.Range("D22").Value = "'5.902"
$ws.Range("E22").Value = "  -0.64%  "
$ws.Range("D23").Value = "28.122.61"
$ws.Range("E23").Value = "  -0.92%  "
$ws.Range("D24").Value = "'11.19"
$ws.Range("E24").Value = "  -0.96%  "
$ws.Range("D25").Value = "'2.090"
$ws.Range("E25").Value = "  -2.14%  "
$ws.Range("D26").Value = "'154.89"
$ws.Range("E26").Value = "  -2.19%  "
$ws.Range("D27").Value = "'20.27"
$ws.Range("E27").Value = "  -1.71%  "
$ws.Range("D28").Value = "1.991.59"
$ws.Range("E28").Value = "  -0.85%  "
$ws.Range("D29").Value = "'2.292"
$ws.Range("E29").Value = "  -3.97%  "
$ws.Range("D30").Value = "'120.78"
$ws.Range("E30").Value = "  -1.69%  "
$ws.Range("D31").Value = "'1.122"
$ws.Range("E31").Value = "  +0.84%  "
$ws.Range("E32").Value = "  +3.69%  "
$ws.Range("E33").Value = "  -1.30%  "
$ws.Range("D34").Value = "'5.557"
$ws.Range("E34").Value = "  -2.88%  "
$ws.Range("E35").Value = "  -2.03%  "
$ws.Range("D36").Value = "'0.06510"
$ws.Range("E36").Value = "  +2.28%  "
$ws.Range("D37").Value = "'0.02294"
$ws.Range("E37").Value = "  -0.82%  "
$ws.Range("D38").Value = "'5.024"
$ws.Range("E38").Value = "  -2.20%  "
$ws.Range("D39").Value = "'8.460"
$ws.Range("E39").Value = "  -3.39%  "
$ws.Range("D40").Value = "'1.453"
$ws.Range("E40").Value = "  +4.70%  "
$ws.Range("D41").Value = "'0.6174"
$ws.Range("E41").Value = "  -3.11%  "
$ws.Range("D42").Value = "'11.08"
$ws.Range("E42").Value = "  -4.20%  "
$ws.Range("D43").Value = "'1.171"
$ws.Range("E43").Value = "  +1.51%  "
$ws.Range("D44").Value = "'0.9998"
$ws.Range("E44").Value = "  -0.02%  "
$ws.Range("D45").Value = "'13.36"
$ws.Range("E45").Value = "  -1.74%  "
$ws.Range("D46").Value = "'3.674"
$ws.Range("E46").Value = "  -0.01%  "
$ws.Range("D47").Value = "'0.5787"
$ws.Range("E47").Value = "  -2.61%  "
$ws.Range("D48").Value = "'125.78"
$ws.Range("E48").Value = "  +0.00%  "
$ws.Range("D49").Value = "'1.191"
$ws.Range("E49").Value = "  +4.07%  "
$ws.Range("D50").Value = "'1.920"
$ws.Range("E50").Value = "  -2.47%  "
$ws.Range("D51").Value = "'0.06810"
$ws.Range("E51").Value = "  -1.32%  "
